$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D2:E51 is treated as Text so numeric-looking strings (e.g. "1.00", "0.999")
# are not auto-converted to numbers, then restore the original (Normal) style afterwards
# so no stray style indices are introduced.
$numRng = $ws.Range("D2:E51")
$numRng.NumberFormat = "@"

$ws.Range("D2").Value = "66.417.82"
$ws.Range("E2").Value = "  -0.61%  "
$ws.Range("D3").Value = "3.316.77"
$ws.Range("E3").Value = "  -1.94%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "190.43"
$ws.Range("E5").Value = "  +2.96%  "
$ws.Range("D6").Value = "562.30"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "0.591"
$ws.Range("E8").Value = "  -1.78%  "
$ws.Range("D9").Value = "3.308.60"
$ws.Range("E9").Value = "  -1.91%  "
$ws.Range("D10").Value = "0.187"
$ws.Range("E10").Value = "  -1.66%  "
$ws.Range("D11").Value = "0.590"
$ws.Range("E11").Value = "  -1.52%  "
$ws.Range("D12").Value = "48.00"
$ws.Range("E12").Value = "  -0.92%  "
$ws.Range("D13").Value = "0.0000272"
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("D14").Value = "8.72"
$ws.Range("E14").Value = "  -0.85%  "
$ws.Range("D15").Value = "3.846.15"
$ws.Range("E15").Value = "  -1.98%  "
$ws.Range("D16").Value = "615.73"
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "18.15"
$ws.Range("E17").Value = "  -1.01%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "66.394.70"
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("E19").Value = "  -0.36%  "
$ws.Range("D20").Value = "3.303.68"
$ws.Range("E20").Value = "  -2.63%  "
$ws.Range("D21").Value = "11.15"
$ws.Range("E21").Value = "  -4.34%  "
$ws.Range("D22").Value = "0.915"
$ws.Range("E22").Value = "  -0.80%  "
$ws.Range("D23").Value = "18.39"
$ws.Range("E23").Value = "  +7.27%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "5.15"
$ws.Range("E24").Value = "  -1.32%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "102.40"
$ws.Range("E25").Value = "  +3.22%  "
$ws.Range("D26").Value = "4.02"
$ws.Range("E26").Value = "  -2.35%  "
$ws.Range("D27").Value = "6.01"
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("D28").Value = "2.76"
$ws.Range("E28").Value = "  +0.66%  "
$ws.Range("E29").Value = "  +2.22%  "
$ws.Range("D30").Value = "8.66"
$ws.Range("E30").Value = "  -2.50%  "
$ws.Range("D31").Value = "30.40"
$ws.Range("E31").Value = "  -1.91%  "
$ws.Range("D32").Value = "4.15"
$ws.Range("E32").Value = "  +5.52%  "
$ws.Range("D33").Value = "6.73"
$ws.Range("E33").Value = "  +5.38%  "
$ws.Range("D34").Value = "569.00"
$ws.Range("E34").Value = "  +2.69%  "
$ws.Range("D35").Value = "11.15"
$ws.Range("E35").Value = "  -0.81%  "
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("D37").Value = "3.758.04"
$ws.Range("E37").Value = "  -3.56%  "
$ws.Range("D38").Value = "57.44"
$ws.Range("E38").Value = "  -1.81%  "
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("B40").Value = "CoreDAO"
$ws.Range("C40").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D40").Value = "3.57"
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "0.0₃0736"
$ws.Range("E41").Value = "  +1.29%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "3.34"
$ws.Range("E42").Value = "  -4.10%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "34.22"
$ws.Range("E43").Value = "  +4.96%  "
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").Value = "0.131"
$ws.Range("E44").Value = "  +1.15%  "
$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").Value = "2.74"
$ws.Range("E45").Value = "  +0.50%  "
$ws.Range("D46").Value = "0.342"
$ws.Range("E46").Value = "  -2.94%  "
$ws.Range("D47").Value = "0.0428"
$ws.Range("E47").Value = "  +1.25%  "
$ws.Range("D48").Value = "3.27"
$ws.Range("E48").Value = "  -0.19%  "
$ws.Range("D49").Value = "0.130"
$ws.Range("E49").Value = "  -1.44%  "
$ws.Range("D50").Value = "2.60"
$ws.Range("E50").Value = "  -4.07%  "
$ws.Range("D51").Value = "0.998"
$ws.Range("E51").Value = "  -0.09%  "

# Restore default styling on the numeric/price/volume columns
$numRng.Style = "Normal"

